$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header/data between column A and column B:
# Column A becomes "Director" data, Column B becomes "Movie" data.
$ws.Range("A1").Value = "Director"
$ws.Range("B1").Value = "Movie"

$ws.Range("A2").Value = "taika waititi"
$ws.Range("B2").Value = "Thor: Love and Thunder"

# Add the new rows of data
$ws.Range("A3").Value = "steven speilberg"
$ws.Range("B3").Value = "West Side Story"
$ws.Range("C3").Value = 76016171

$ws.Range("A4").Value = "james cameron"
$ws.Range("B4").Value = "Duets"
$ws.Range("C4").Value = 6620242

$ws.Range("A5").Value = "christopher nolan"
$ws.Range("B5").Value = "Quay"
$ws.Range("C5").Value = 51858
